$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) and Volume(1h) (E) columns for the rows being
# updated so the numeric-looking strings are stored as text, matching the original
# inline-string cell type instead of being auto-coerced to numbers.
$cells = @{
    'D2' = '286.89'
    'E2' = '1.51%'
    'D3' = '29.43'
    'E3' = '4.07%'
    'D4' = '5.071'
    'E4' = '0.72%'
    'D5' = '0.06719'
    'E5' = '2.96%'
    'D6' = '7.330'
    'E6' = '1.24%'
    'D7' = '1.374'
    'E7' = '-1.46%'
    'D8' = '0.9027'
    'E8' = '-1.73%'
    'D9' = '0.1594'
    'E9' = '2.86%'
    'D10' = '0.06959'
    'E10' = '7.30%'
    'D11' = '0.07630'
    'E11' = '1.14%'
    'D12' = '0.02927'
    'E12' = '6.76%'
    'D13' = '0.08989'
    'E13' = '0.33%'
    'D14' = '0.001567'
    'E14' = '-0.93%'
    'D15' = '0.04481'
    'E15' = '1.67%'
    'D16' = '0.0006481'
    'E16' = '1.17%'
    'D17' = '0.006188'
    'E17' = '1.34%'
    'D18' = '3.452'
    'E18' = '0.21%'
    'D19' = '3.452'
    'E19' = '2.53%'
    'D20' = '2.231'
    'E20' = '-0.12%'
    'D21' = '0.3204'
    'E21' = '0.53%'
    'D22' = '0.1314'
    'E22' = '2.13%'
    'D23' = '3.945'
    'E23' = '-0.70%'
    'D24' = '0.1561'
    'E24' = '1.10%'
    'D25' = '0.001201'
    'E25' = '1.80%'
    'D26' = '0.004366'
    'E26' = '-1.13%'
    'E27' = '0.15%'
    'D28' = '0.0001619'
    'E28' = '0.11%'
    'D40' = '0.04248'
    'E40' = '2.95%'
    'D41' = '0.006802'
    'E41' = '1.37%'
    'D42' = '0.1238'
    'E42' = '0.87%'
    'E43' = '4.36%'
    'D44' = '0.01261'
    'E44' = '4.60%'
    'D45' = '0.00005750'
    'E45' = '3.18%'
    'D46' = '1.962'
    'E46' = '-0.18%'
    'D47' = '0.01501'
    'E47' = '-18.91%'
}

foreach ($ref in $cells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$ref]
}

Write-Output "Updated $($cells.Count) cells"
